$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format so the numeric-looking strings are not
# auto-converted to numbers by Excel, then clear the format delta back to
# the original (unstyled) state once the text values are in place.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.114.89"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "1.891.99"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "306.66"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.5157"
$ws.Range("E7").Value = "  +2.12%  "
$ws.Range("D8").Value = "0.3762"
$ws.Range("E8").Value = "  +3.21%  "
$ws.Range("D9").Value = "0.07199"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "21.16"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").Value = "0.9024"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "0.07660"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").Value = "1.868.57"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "94.41"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "5.248"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "0.000008501"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "14.45"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "27.138.14"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "5.062"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "2.129.14"
$ws.Range("E22").Value = "  +2.55%  "
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "2.284"
$ws.Range("E25").Value = "  +9.93%  "
$ws.Range("D26").Value = "146.60"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "1.740"
$ws.Range("E27").Value = "  -2.84%  "
$ws.Range("D28").Value = "18.08"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").Value = "114.51"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").Value = "4.934"
$ws.Range("E30").Value = "  +5.56%  "
$ws.Range("D31").Value = "4.802"
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("D32").Value = "0.09209"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "0.05051"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "1.236"
$ws.Range("E34").Value = "  +7.45%  "
$ws.Range("D35").Value = "0.7686"
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("D36").Value = "2.989"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "3.281"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "2.595"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("D39").Value = "0.5610"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "0.01991"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").Value = "1.072"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "9.079"
$ws.Range("E42").Value = "  +6.59%  "
$ws.Range("D43").Value = "6.639"
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("D44").Value = "118.60"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "0.1503"
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("D46").Value = "0.4820"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "10.13"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").Value = "37.62"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").Value = "64.05"
$ws.Range("E51").Value = "  +1.89%  "

$ws.Range("D2:E51").ClearFormats()
